$wb = $excel.ActiveWorkbook

# Helper: write a value as literal TEXT (never auto-coerced to a number/date),
# matching how the source data stores numeric-looking strings such as fund
# codes ("009630") or percentages ("7.35") as plain text. The leading "'"
# forces Excel to keep the literal text; ClearFormats() then drops the
# resulting quote-prefix formatting flag so the cell keeps the worksheet's
# default (unstyled) look, same as the other text cells around it.
function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
    $cell.ClearFormats()
}

# The workbook currently ends with a "总计" (grand-total) summary sheet.
# A new quarter of data (2022-Q1) has arrived, so:
#  1. the existing "总计" sheet is repurposed into the new "2022-Q1" detail
#     sheet (same 8-column fund-holdings layout as the other quarter sheets);
#  2. a fresh "总计" sheet is appended after it, with the running summary
#     table updated to include the new 2022-Q1 row.
$oldTotal = $wb.Worksheets.Item("总计")
$oldTotal.Name = "2022-Q1"

$newTotal = $wb.Worksheets.Add($null, $oldTotal)
$newTotal.Name = "总计"

$q1 = $wb.Worksheets.Item("2022-Q1")
$total = $wb.Worksheets.Item("总计")

# Style template: an existing header cell that already carries the workbook's
# standard bold / bordered / centered formatting (the style used for the
# header row and index column A on every quarter/total sheet). Copying it
# keeps new cells visually consistent with the rest of the workbook.
$styleTemplate = $wb.Worksheets.Item("2021-Q4").Cells.Item(1, 2)

# --- "2022-Q1" sheet: fund holding detail table --------------------------
$styleTemplate.Copy($q1.Cells.Item(1, 2))
$q1.Cells.Item(1, 2).Value = "基金代码"
$styleTemplate.Copy($q1.Cells.Item(1, 3))
$q1.Cells.Item(1, 3).Value = "基金名称"
$styleTemplate.Copy($q1.Cells.Item(1, 4))
$q1.Cells.Item(1, 4).Value = "基金规模"
$styleTemplate.Copy($q1.Cells.Item(1, 5))
$q1.Cells.Item(1, 5).Value = "股票总仓位"
$styleTemplate.Copy($q1.Cells.Item(1, 6))
$q1.Cells.Item(1, 6).Value = "仓位占比"
$styleTemplate.Copy($q1.Cells.Item(1, 7))
$q1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$styleTemplate.Copy($q1.Cells.Item(1, 8))
$q1.Cells.Item(1, 8).Value = "仓位排名"

$styleTemplate.Copy($q1.Cells.Item(2, 1))
$q1.Cells.Item(2, 1).Value = 0
Set-TextValue $q1.Cells.Item(2, 2) "009630"
Set-TextValue $q1.Cells.Item(2, 3) "浦银安盛ESG责任投资混合A"
Set-TextValue $q1.Cells.Item(2, 4) "15.61"
Set-TextValue $q1.Cells.Item(2, 5) "80.10"
Set-TextValue $q1.Cells.Item(2, 6) "7.35"
Set-TextValue $q1.Cells.Item(2, 7) "1.1473"
$q1.Cells.Item(2, 8).Value = 2

$styleTemplate.Copy($q1.Cells.Item(3, 1))
$q1.Cells.Item(3, 1).Value = 1
Set-TextValue $q1.Cells.Item(3, 2) "002666"
Set-TextValue $q1.Cells.Item(3, 3) "前海开源沪港深创新成长灵活配置混合A"
Set-TextValue $q1.Cells.Item(3, 4) "11.96"
Set-TextValue $q1.Cells.Item(3, 5) "81.64"
Set-TextValue $q1.Cells.Item(3, 6) "6.85"
Set-TextValue $q1.Cells.Item(3, 7) "0.8193"
$q1.Cells.Item(3, 8).Value = 5

$styleTemplate.Copy($q1.Cells.Item(4, 1))
$q1.Cells.Item(4, 1).Value = 2
Set-TextValue $q1.Cells.Item(4, 2) "001230"
Set-TextValue $q1.Cells.Item(4, 3) "鹏华医药科技股票"
Set-TextValue $q1.Cells.Item(4, 4) "16.65"
Set-TextValue $q1.Cells.Item(4, 5) "81.88"
Set-TextValue $q1.Cells.Item(4, 6) "3.81"
Set-TextValue $q1.Cells.Item(4, 7) "0.6344"
$q1.Cells.Item(4, 8).Value = 8

$styleTemplate.Copy($q1.Cells.Item(5, 1))
$q1.Cells.Item(5, 1).Value = 3
Set-TextValue $q1.Cells.Item(5, 2) "001766"
Set-TextValue $q1.Cells.Item(5, 3) "上投摩根医疗健康股票"
Set-TextValue $q1.Cells.Item(5, 4) "10.35"
Set-TextValue $q1.Cells.Item(5, 5) "80.54"
Set-TextValue $q1.Cells.Item(5, 6) "5.62"
Set-TextValue $q1.Cells.Item(5, 7) "0.5817"
$q1.Cells.Item(5, 8).Value = 5

$styleTemplate.Copy($q1.Cells.Item(6, 1))
$q1.Cells.Item(6, 1).Value = 4
Set-TextValue $q1.Cells.Item(6, 2) "009631"
Set-TextValue $q1.Cells.Item(6, 3) "浦银安盛ESG责任投资混合C"
Set-TextValue $q1.Cells.Item(6, 4) "5.74"
Set-TextValue $q1.Cells.Item(6, 5) "80.10"
Set-TextValue $q1.Cells.Item(6, 6) "7.35"
Set-TextValue $q1.Cells.Item(6, 7) "0.4219"
$q1.Cells.Item(6, 8).Value = 2

$styleTemplate.Copy($q1.Cells.Item(7, 1))
$q1.Cells.Item(7, 1).Value = 5
Set-TextValue $q1.Cells.Item(7, 2) "007066"
Set-TextValue $q1.Cells.Item(7, 3) "浦银安盛先进制造混合A"
Set-TextValue $q1.Cells.Item(7, 4) "3.66"
Set-TextValue $q1.Cells.Item(7, 5) "74.55"
Set-TextValue $q1.Cells.Item(7, 6) "7.56"
Set-TextValue $q1.Cells.Item(7, 7) "0.2767"
$q1.Cells.Item(7, 8).Value = 2

$styleTemplate.Copy($q1.Cells.Item(8, 1))
$q1.Cells.Item(8, 1).Value = 6
Set-TextValue $q1.Cells.Item(8, 2) "002667"
Set-TextValue $q1.Cells.Item(8, 3) "前海开源沪港深创新成长灵活配置混合C"
Set-TextValue $q1.Cells.Item(8, 4) "3.25"
Set-TextValue $q1.Cells.Item(8, 5) "81.64"
Set-TextValue $q1.Cells.Item(8, 6) "6.85"
Set-TextValue $q1.Cells.Item(8, 7) "0.2226"
$q1.Cells.Item(8, 8).Value = 5

$styleTemplate.Copy($q1.Cells.Item(9, 1))
$q1.Cells.Item(9, 1).Value = 7
Set-TextValue $q1.Cells.Item(9, 2) "007067"
Set-TextValue $q1.Cells.Item(9, 3) "浦银安盛先进制造混合C"
Set-TextValue $q1.Cells.Item(9, 4) "2.27"
Set-TextValue $q1.Cells.Item(9, 5) "74.55"
Set-TextValue $q1.Cells.Item(9, 6) "7.56"
Set-TextValue $q1.Cells.Item(9, 7) "0.1716"
$q1.Cells.Item(9, 8).Value = 2

$styleTemplate.Copy($q1.Cells.Item(10, 1))
$q1.Cells.Item(10, 1).Value = 8
Set-TextValue $q1.Cells.Item(10, 2) "003581"
Set-TextValue $q1.Cells.Item(10, 3) "新疆前海联合国民健康产业灵活配置混合A"
Set-TextValue $q1.Cells.Item(10, 4) "2.66"
Set-TextValue $q1.Cells.Item(10, 5) "94.65"
Set-TextValue $q1.Cells.Item(10, 6) "4.64"
Set-TextValue $q1.Cells.Item(10, 7) "0.1234"
$q1.Cells.Item(10, 8).Value = 9

$styleTemplate.Copy($q1.Cells.Item(11, 1))
$q1.Cells.Item(11, 1).Value = 9
Set-TextValue $q1.Cells.Item(11, 2) "012027"
Set-TextValue $q1.Cells.Item(11, 3) "光大保德信安阳一年持有期混合型证券投资基金A"
Set-TextValue $q1.Cells.Item(11, 4) "15.22"
Set-TextValue $q1.Cells.Item(11, 5) "22.05"
Set-TextValue $q1.Cells.Item(11, 6) "0.75"
Set-TextValue $q1.Cells.Item(11, 7) "0.1142"
$q1.Cells.Item(11, 8).Value = 9

$styleTemplate.Copy($q1.Cells.Item(12, 1))
$q1.Cells.Item(12, 1).Value = 10
Set-TextValue $q1.Cells.Item(12, 2) "001415"
Set-TextValue $q1.Cells.Item(12, 3) "信诚新锐回报灵活配置混合A"
Set-TextValue $q1.Cells.Item(12, 4) "9.07"
Set-TextValue $q1.Cells.Item(12, 5) "24.72"
Set-TextValue $q1.Cells.Item(12, 6) "0.66"
Set-TextValue $q1.Cells.Item(12, 7) "0.0599"
$q1.Cells.Item(12, 8).Value = 7

$styleTemplate.Copy($q1.Cells.Item(13, 1))
$q1.Cells.Item(13, 1).Value = 11
Set-TextValue $q1.Cells.Item(13, 2) "012028"
Set-TextValue $q1.Cells.Item(13, 3) "光大保德信安阳一年持有期混合型证券投资基金C"
Set-TextValue $q1.Cells.Item(13, 4) "7.68"
Set-TextValue $q1.Cells.Item(13, 5) "22.05"
Set-TextValue $q1.Cells.Item(13, 6) "0.75"
Set-TextValue $q1.Cells.Item(13, 7) "0.0576"
$q1.Cells.Item(13, 8).Value = 9

$styleTemplate.Copy($q1.Cells.Item(14, 1))
$q1.Cells.Item(14, 1).Value = 12
Set-TextValue $q1.Cells.Item(14, 2) "003234"
Set-TextValue $q1.Cells.Item(14, 3) "信诚至利灵活配置混合A"
Set-TextValue $q1.Cells.Item(14, 4) "8.99"
Set-TextValue $q1.Cells.Item(14, 5) "22.05"
Set-TextValue $q1.Cells.Item(14, 6) "0.61"
Set-TextValue $q1.Cells.Item(14, 7) "0.0548"
$q1.Cells.Item(14, 8).Value = 5

$styleTemplate.Copy($q1.Cells.Item(15, 1))
$q1.Cells.Item(15, 1).Value = 13
Set-TextValue $q1.Cells.Item(15, 2) "001402"
Set-TextValue $q1.Cells.Item(15, 3) "信诚新选回报灵活配置混合A"
Set-TextValue $q1.Cells.Item(15, 4) "8.37"
Set-TextValue $q1.Cells.Item(15, 5) "22.05"
Set-TextValue $q1.Cells.Item(15, 6) "0.61"
Set-TextValue $q1.Cells.Item(15, 7) "0.0511"
$q1.Cells.Item(15, 8).Value = 7

$styleTemplate.Copy($q1.Cells.Item(16, 1))
$q1.Cells.Item(16, 1).Value = 14
Set-TextValue $q1.Cells.Item(16, 2) "004157"
Set-TextValue $q1.Cells.Item(16, 3) "信诚至诚灵活配置混合A"
Set-TextValue $q1.Cells.Item(16, 4) "7.32"
Set-TextValue $q1.Cells.Item(16, 5) "22.71"
Set-TextValue $q1.Cells.Item(16, 6) "0.62"
Set-TextValue $q1.Cells.Item(16, 7) "0.0454"
$q1.Cells.Item(16, 8).Value = 7

$styleTemplate.Copy($q1.Cells.Item(17, 1))
$q1.Cells.Item(17, 1).Value = 15
Set-TextValue $q1.Cells.Item(17, 2) "002046"
Set-TextValue $q1.Cells.Item(17, 3) "信诚新锐回报灵活配置混合B"
Set-TextValue $q1.Cells.Item(17, 4) "5.19"
Set-TextValue $q1.Cells.Item(17, 5) "24.72"
Set-TextValue $q1.Cells.Item(17, 6) "0.66"
Set-TextValue $q1.Cells.Item(17, 7) "0.0343"
$q1.Cells.Item(17, 8).Value = 7

$styleTemplate.Copy($q1.Cells.Item(18, 1))
$q1.Cells.Item(18, 1).Value = 16
Set-TextValue $q1.Cells.Item(18, 2) "003235"
Set-TextValue $q1.Cells.Item(18, 3) "信诚至利灵活配置混合C"
Set-TextValue $q1.Cells.Item(18, 4) "5.30"
Set-TextValue $q1.Cells.Item(18, 5) "22.05"
Set-TextValue $q1.Cells.Item(18, 6) "0.61"
Set-TextValue $q1.Cells.Item(18, 7) "0.0323"
$q1.Cells.Item(18, 8).Value = 5

$styleTemplate.Copy($q1.Cells.Item(19, 1))
$q1.Cells.Item(19, 1).Value = 17
Set-TextValue $q1.Cells.Item(19, 2) "010703"
Set-TextValue $q1.Cells.Item(19, 3) "财通智选消费股票A"
Set-TextValue $q1.Cells.Item(19, 4) "1.12"
Set-TextValue $q1.Cells.Item(19, 5) "92.43"
Set-TextValue $q1.Cells.Item(19, 6) "2.85"
Set-TextValue $q1.Cells.Item(19, 7) "0.0319"
$q1.Cells.Item(19, 8).Value = 5

$styleTemplate.Copy($q1.Cells.Item(20, 1))
$q1.Cells.Item(20, 1).Value = 18
Set-TextValue $q1.Cells.Item(20, 2) "002030"
Set-TextValue $q1.Cells.Item(20, 3) "信诚新选回报灵活配置混合B"
Set-TextValue $q1.Cells.Item(20, 4) "3.85"
Set-TextValue $q1.Cells.Item(20, 5) "22.05"
Set-TextValue $q1.Cells.Item(20, 6) "0.61"
Set-TextValue $q1.Cells.Item(20, 7) "0.0235"
$q1.Cells.Item(20, 8).Value = 7

$styleTemplate.Copy($q1.Cells.Item(21, 1))
$q1.Cells.Item(21, 1).Value = 19
Set-TextValue $q1.Cells.Item(21, 2) "010704"
Set-TextValue $q1.Cells.Item(21, 3) "财通智选消费股票C"
Set-TextValue $q1.Cells.Item(21, 4) "0.54"
Set-TextValue $q1.Cells.Item(21, 5) "92.43"
Set-TextValue $q1.Cells.Item(21, 6) "2.85"
Set-TextValue $q1.Cells.Item(21, 7) "0.0154"
$q1.Cells.Item(21, 8).Value = 5

$styleTemplate.Copy($q1.Cells.Item(22, 1))
$q1.Cells.Item(22, 1).Value = 20
Set-TextValue $q1.Cells.Item(22, 2) "007111"
Set-TextValue $q1.Cells.Item(22, 3) "新疆前海联合国民健康产业灵活配置混合C"
Set-TextValue $q1.Cells.Item(22, 4) "0.32"
Set-TextValue $q1.Cells.Item(22, 5) "94.65"
Set-TextValue $q1.Cells.Item(22, 6) "4.64"
Set-TextValue $q1.Cells.Item(22, 7) "0.0148"
$q1.Cells.Item(22, 8).Value = 9

$styleTemplate.Copy($q1.Cells.Item(23, 1))
$q1.Cells.Item(23, 1).Value = 21
Set-TextValue $q1.Cells.Item(23, 2) "004158"
Set-TextValue $q1.Cells.Item(23, 3) "信诚至诚灵活配置混合B"
Set-TextValue $q1.Cells.Item(23, 4) "2.18"
Set-TextValue $q1.Cells.Item(23, 5) "22.71"
Set-TextValue $q1.Cells.Item(23, 6) "0.62"
Set-TextValue $q1.Cells.Item(23, 7) "0.0135"
$q1.Cells.Item(23, 8).Value = 7

# --- "总计" sheet: quarter-over-quarter summary table --------------------
$styleTemplate.Copy($total.Cells.Item(1, 2))
$total.Cells.Item(1, 2).Value = "日期"
$styleTemplate.Copy($total.Cells.Item(1, 3))
$total.Cells.Item(1, 3).Value = "持有数量(只)"
$styleTemplate.Copy($total.Cells.Item(1, 4))
$total.Cells.Item(1, 4).Value = "持有市值(亿元)"

$styleTemplate.Copy($total.Cells.Item(2, 1))
$total.Cells.Item(2, 1).Value = 0
Set-TextValue $total.Cells.Item(2, 2) "2022-Q1"
$total.Cells.Item(2, 3).Value = 22
$total.Cells.Item(2, 4).Value = 4.95

$styleTemplate.Copy($total.Cells.Item(3, 1))
$total.Cells.Item(3, 1).Value = 1
Set-TextValue $total.Cells.Item(3, 2) "2021-Q4"
$total.Cells.Item(3, 3).Value = 9
$total.Cells.Item(3, 4).Value = 5.66

$styleTemplate.Copy($total.Cells.Item(4, 1))
$total.Cells.Item(4, 1).Value = 2
Set-TextValue $total.Cells.Item(4, 2) "2021-Q3"
$total.Cells.Item(4, 3).Value = 19
$total.Cells.Item(4, 4).Value = 5.82

$styleTemplate.Copy($total.Cells.Item(5, 1))
$total.Cells.Item(5, 1).Value = 3
Set-TextValue $total.Cells.Item(5, 2) "2021-Q2"
$total.Cells.Item(5, 3).Value = 11
$total.Cells.Item(5, 4).Value = 19.93

$styleTemplate.Copy($total.Cells.Item(6, 1))
$total.Cells.Item(6, 1).Value = 4
Set-TextValue $total.Cells.Item(6, 2) "2021-Q1"
$total.Cells.Item(6, 3).Value = 12
$total.Cells.Item(6, 4).Value = 19.07

$styleTemplate.Copy($total.Cells.Item(7, 1))
$total.Cells.Item(7, 1).Value = 5
Set-TextValue $total.Cells.Item(7, 2) "2020-Q4"
$total.Cells.Item(7, 3).Value = 23
$total.Cells.Item(7, 4).Value = 25.6

